# Commit: "added EnrollCourseTest,Cart_page,EnrollPage,added few methods in HomePage"
#
# The xlsx-relevant portion of this commit populates the previously-empty
# "courses_details" sheet with the sample data used by the new
# EnrollCourseTest (a course name, an address and a phone number), styles
# the header row with the workbook's existing yellow highlight fill, and
# makes "courses_details" the active/selected sheet (it was "new_users"
# before).

$wb = $excel.ActiveWorkbook

# courses_details is the 3rd tab: new_users, login_details, courses_details, ...
$ws = $wb.Worksheets.Item(3)

# Header row
$ws.Range("A1").Value = "Course"
$ws.Range("B1").Value = "Address"
$ws.Range("C1").Value = "Phone"

# Highlight the header row with the workbook's existing yellow fill
# (same RGB FFFF00 already used elsewhere in the workbook for header rows).
$ws.Range("A1:C1").Interior.Color = 65535

# Data row
$ws.Range("A2").Value = "Selenium"
$ws.Range("B2").Value = "Hyderabad"
$ws.Range("C2").Value = 1234567891

# Make this the active sheet / tab, with A2 selected, matching the new
# workbook-level activeTab and the sheet's tabSelected + selection.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
